# "Generate Report for Handoff"
# The localization-status report is regenerated: the row for
# 9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md now sorts first (still
# "Handed back: in sync with en-US"), and the row for
# 34231130-4ec4-4dc4-a8d4-fa28d37ba982.md moves to the second slot and
# is now "Ready for handoff" with a refreshed handoff datetime and an
# error detail noting the handback file is stale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$ov.Range("B2").Value = "e2e\9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"

$ov.Range("A3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$ov.Range("B3").Value = "e2e\34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-20 12:51:00"

$ovLinks = @($ov.Hyperlinks)
$ovLinks[0].TextToDisplay = "e2e\9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$ovLinks[1].TextToDisplay = "e2e\34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$zh.Range("G2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.bacbd0f75d3e3dd13285c18d560b44f1d33c5a65.zh-cn.xlf"
$zh.Range("I2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$zh.Range("J2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.bacbd0f75d3e3dd13285c18d560b44f1d33c5a65.zh-cn.xlf"

$zh.Range("A3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.e60f937d3ba5d8d205564b504a553e95fb736470.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-20 12:50:56"
$zh.Range("I3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$zh.Range("J3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.e60f937d3ba5d8d205564b504a553e95fb736470.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6c455ded9406a780f9ac6869ab8d703b684720f/e2e/34231130-4ec4-4dc4-a8d4-fa28d37ba982.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3204c8756e1606d9524e4afc910618e6895175c6/e2e/34231130-4ec4-4dc4-a8d4-fa28d37ba982.md."

$zhLinks = @($zh.Hyperlinks)
$zhLinks[0].TextToDisplay = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$zhLinks[1].TextToDisplay = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$zhLinks[2].TextToDisplay = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$zhLinks[3].TextToDisplay = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"

$zh.Columns.Item(16).ColumnWidth = 39.1667

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$de.Range("G2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.bacbd0f75d3e3dd13285c18d560b44f1d33c5a65.de-de.xlf"
$de.Range("I2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$de.Range("J2").Value = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.bacbd0f75d3e3dd13285c18d560b44f1d33c5a65.de-de.xlf"

$de.Range("A3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.e60f937d3ba5d8d205564b504a553e95fb736470.de-de.xlf"
$de.Range("H3").Value = "2016-08-20 12:51:00"
$de.Range("I3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$de.Range("J3").Value = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.e60f937d3ba5d8d205564b504a553e95fb736470.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a6c455ded9406a780f9ac6869ab8d703b684720f/e2e/34231130-4ec4-4dc4-a8d4-fa28d37ba982.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3204c8756e1606d9524e4afc910618e6895175c6/e2e/34231130-4ec4-4dc4-a8d4-fa28d37ba982.md."

$deLinks = @($de.Hyperlinks)
$deLinks[0].TextToDisplay = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$deLinks[1].TextToDisplay = "9f30a533-d9f1-4659-bd9d-2033d0eb0ef6.md"
$deLinks[2].TextToDisplay = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"
$deLinks[3].TextToDisplay = "34231130-4ec4-4dc4-a8d4-fa28d37ba982.md"

$de.Columns.Item(16).ColumnWidth = 39.1667
